$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# This string is shared across the Overview summary columns and each
# locale sheet's Status column, so every occurrence is rewritten.

# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: column C (Status), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: column C (Status), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Re-fit the Status-related columns to the new, shorter text ---
# "In Translation" is shorter than "Ready for handoff", so the columns
# that were sized to that text shrink accordingly (Overview!E:F and the
# Status column on each locale sheet).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
